$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.184.49"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.69%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.580.72"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.06%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.02"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.78%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.48"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.14%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.31%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.22%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.590.59"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.55%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.94%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.59%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.160"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +11.83%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.346"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.50%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.033.47"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.13%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.151.49"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.70%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.92"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.61%  "

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.66%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.586.64"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.24%  "

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.49%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.83"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.81%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.37"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.20%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.39"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.72%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.19%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.17"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.52%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.462"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.00%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.55%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.99%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.31"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.13%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.57%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.10%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "161.30"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.21%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.28%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.63%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.92"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.32%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.72%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.01%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.873"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.13%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.45"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.07%  "

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.56%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "294.18"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.52%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.65"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.22%  "

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.42%  "

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "131.90"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.89%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0972"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.43%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.50%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0534"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.60%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.12%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.05"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.52%  "

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.97%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.52"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.40%  "
